# Powerpoint writer: consolidate text run nodes.
# Merge adjacent runs that share identical formatting by rewriting the
# combined character span's Text, which causes the host to fold the
# span into a single run instead of leaving the original per-word runs.

$p = $ppt.ActivePresentation

# --- Slide 1: Title "Header" / " " / "with" / " " / "inline code" ---
# -> "Header " / "with " / "inline code"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 7).Text = "Header "   # "Header" + " "
$tr1.Characters(8, 5).Text = "with "     # "with" + " "

# --- Slide 2: Title "Syntax" / " " / "highlighting" ---
# -> "Syntax " / "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 7).Text = "Syntax "   # "Syntax" + " "

# --- Slide 3: Title "Two" / " " / "column" / " " / "slide" ---
# -> "Two " / "column " / "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 4).Text = "Two "      # "Two" + " "
$tr3.Characters(5, 7).Text = "column "   # "column" + " "
